$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (10th column, "Priority"), shifting
# everything from J onward to the right by one. Excel carries the left
# neighbor's formatting into the new column automatically, matching the
# surrounding "Status"/"Priority" header & body cell styling.
$ws.Columns("J").Insert()

# Header text for the newly inserted "Severidad" column.
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Value = "Severidad"

# Data cell for the new column stays formatted like its row neighbors but
# carries no value.
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122) | Out-Null
$ws.Range("J3").ClearContents()

$excel.CutCopyMode = 0

# Row 3 re-wraps slightly shorter once the new column is present.
$ws.Rows("3").RowHeight = 114.75

# Update the active selection to the new data cell for the column, matching
# the author's final selection state.
$ws.Range("J3").Select()
